$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 12 de Abril de 2020 a las 14:52"

# Update Galicia's row (row 6) statistics
$ws.Range("B6").Value = 7336
$ws.Range("C6").Value = 1143
$ws.Range("D6").Value = 5833
$ws.Range("E6").Value = 360
